$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2023-09-25 Monday" "2023-09-26 Tuesday"

Replace-Text "27÷8=" "41÷6="
Replace-Text "49÷9=" "12÷8="
Replace-Text "78÷6=" "39÷4="
Replace-Text "52÷3=" "12÷2="
Replace-Text "97÷7=" "70÷7="
Replace-Text "21÷7=" "54÷9="
Replace-Text "99÷6=" "20÷2="
Replace-Text "38÷4=" "80÷9="
Replace-Text "98÷7=" "79÷3="
Replace-Text "24÷9=" "99÷9="
Replace-Text "31÷7=" "73÷2="
Replace-Text "50÷8=" "94÷7="
Replace-Text "97÷5=" "72÷8="
Replace-Text "66÷2=" "52÷5="
Replace-Text "51÷5=" "83÷8="
Replace-Text "45÷2=" "79÷3="
Replace-Text "80÷5=" "80÷6="
Replace-Text "96÷4=" "53÷5="
Replace-Text "32÷4=" "38÷4="
Replace-Text "21÷3=" "21÷4="
Replace-Text "34÷4=" "12÷6="
Replace-Text "24÷3=" "52÷8="
Replace-Text "29÷8=" "85÷7="
Replace-Text "84÷9=" "79÷2="
Replace-Text "99÷7=" "31÷5="
